# Commit: "adding last crud operations"
# Adds a boolean "Test Passed" flag (column E on Clients, column G on Projects)
# to every existing data row, marking the CRUD test cases as passed.

$wb = $excel.ActiveWorkbook

# --- Sheet "Clients": set boolean TRUE in column E for rows 2..50 ---
$wsClients = $wb.Worksheets.Item("Clients")
for ($r = 2; $r -le 50; $r++) {
    $wsClients.Cells.Item($r, 5).Value = $true
}
[void]$wsClients.Range("E2").Select()

# --- Sheet "Projects": set boolean TRUE in column G for rows 2..46 ---
$wsProjects = $wb.Worksheets.Item("Projects")
for ($r = 2; $r -le 46; $r++) {
    $wsProjects.Cells.Item($r, 7).Value = $true
}

# Leave "Projects" as the active sheet/tab, selection on G2, matching
# the last state the user left the workbook in after the edits.
[void]$wsProjects.Activate()
[void]$wsProjects.Range("G2").Select()
